# Add a new record (row 15) to the "CredencialesRealizadas" sheet,
# mirroring the existing rows' layout/columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 15

# Column A ("Unnamed: 0") is blank for this record, like most other rows.
# A leading apostrophe forces an explicit (empty) text value instead of
# Excel dropping the write entirely; resetting the style afterwards keeps
# the cell free of the "text quote prefix" formatting flag.
$ws.Cells.Item($row, 1).Value = "'"
$ws.Cells.Item($row, 1).Style = "Normal"

# Folio - stored as text "12345" (not a number) to match the source data.
$ws.Cells.Item($row, 2).Value = "'12345"
$ws.Cells.Item($row, 2).Style = "Normal"

$ws.Cells.Item($row, 3).Value = "Juan"
$ws.Cells.Item($row, 4).Value = "Carlos"
$ws.Cells.Item($row, 5).Value = "Calderon"
$ws.Cells.Item($row, 6).Value = "Davila"
$ws.Cells.Item($row, 7).Value = "Sub-director"
$ws.Cells.Item($row, 8).Value = "Subdirector Administrativo"
$ws.Cells.Item($row, 9).Value = "23/03/2023"

# Vigencia - a date serial, formatted the same way as the row above it.
$ws.Cells.Item($row, 10).Value = 45009
$ws.Cells.Item($row, 10).Style = $ws.Cells.Item($row - 1, 10).Style
$ws.Cells.Item($row, 10).NumberFormat = $ws.Cells.Item($row - 1, 10).NumberFormat

$ws.Cells.Item($row, 11).Value = 140198
$ws.Cells.Item($row, 12).Value = "Niels Henrick Navarrete Manzanilla"
$ws.Cells.Item($row, 13).Value = "C:/Users/MrJua/Downloads/20230318_002304.jpg"

$wb.Save()
